$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D and B/D styling are identical across the new rows (842-853), so set them in bulk first.
$ws.Range("A842:A853").Value = "Entrainement"
$ws.Range("B842:B853").Value = 45974
$ws.Range("B842:B853").NumberFormat = "m/d/yy"
$ws.Range("C842:C853").Value = "Global"
$ws.Range("D842:D853").Value = "J-3"
$ws.Range("D842:D853").HorizontalAlignment = -4108

# Per-row player stats (columns E through V).
# Row 842
$ws.Cells.Item(842, 5).Value = "Mattheo Haon"
$ws.Cells.Item(842, 6).Value = "right back"
$ws.Cells.Item(842, 7).Value = "01:30:43"
$ws.Cells.Item(842, 8).Value = 6.33
$ws.Cells.Item(842, 9).Value = 0.22
$ws.Cells.Item(842, 10).Value = 6.11
$ws.Cells.Item(842, 11).Value = 0.22
$ws.Cells.Item(842, 12).Value = 0.01
$ws.Cells.Item(842, 13).Value = 0.0
$ws.Cells.Item(842, 14).Value = 0.0
$ws.Cells.Item(842, 15).Value = 0.0
$ws.Cells.Item(842, 16).Value = 4.14
$ws.Cells.Item(842, 17).Value = 20.33
$ws.Cells.Item(842, 18).Value = 4.62
$ws.Cells.Item(842, 19).Value = 35.0
$ws.Cells.Item(842, 20).Value = 4.0
$ws.Cells.Item(842, 21).Value = 31.0
$ws.Cells.Item(842, 22).Value = 5.0

# Row 843
$ws.Cells.Item(843, 5).Value = "Naim Ighbane"
$ws.Cells.Item(843, 6).Value = "center back"
$ws.Cells.Item(843, 7).Value = "01:30:13"
$ws.Cells.Item(843, 8).Value = 5.4
$ws.Cells.Item(843, 9).Value = 0.04
$ws.Cells.Item(843, 10).Value = 5.36
$ws.Cells.Item(843, 11).Value = 0.04
$ws.Cells.Item(843, 12).Value = 0.0
$ws.Cells.Item(843, 13).Value = 0.0
$ws.Cells.Item(843, 14).Value = 0.0
$ws.Cells.Item(843, 15).Value = 0.0
$ws.Cells.Item(843, 16).Value = 2.83
$ws.Cells.Item(843, 17).Value = 18.09
$ws.Cells.Item(843, 18).Value = 3.37
$ws.Cells.Item(843, 19).Value = 8.0
$ws.Cells.Item(843, 20).Value = 0.0
$ws.Cells.Item(843, 21).Value = 6.0
$ws.Cells.Item(843, 22).Value = 2.0

# Row 844
$ws.Cells.Item(844, 5).Value = "Sofiane Belle"
$ws.Cells.Item(844, 6).Value = "left forward"
$ws.Cells.Item(844, 7).Value = "01:27:30"
$ws.Cells.Item(844, 8).Value = 5.15
$ws.Cells.Item(844, 9).Value = 0.27
$ws.Cells.Item(844, 10).Value = 4.87
$ws.Cells.Item(844, 11).Value = 0.25
$ws.Cells.Item(844, 12).Value = 0.02
$ws.Cells.Item(844, 13).Value = 0.0
$ws.Cells.Item(844, 14).Value = 0.0
$ws.Cells.Item(844, 15).Value = 0.0
$ws.Cells.Item(844, 16).Value = 3.36
$ws.Cells.Item(844, 17).Value = 23.84
$ws.Cells.Item(844, 18).Value = 4.01
$ws.Cells.Item(844, 19).Value = 15.0
$ws.Cells.Item(844, 20).Value = 1.0
$ws.Cells.Item(844, 21).Value = 12.0
$ws.Cells.Item(844, 22).Value = 1.0

# Row 845
$ws.Cells.Item(845, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(845, 6).Value = "center midfield"
$ws.Cells.Item(845, 7).Value = "01:29:34"
$ws.Cells.Item(845, 8).Value = 7.38
$ws.Cells.Item(845, 9).Value = 0.44
$ws.Cells.Item(845, 10).Value = 6.93
$ws.Cells.Item(845, 11).Value = 0.36
$ws.Cells.Item(845, 12).Value = 0.1
$ws.Cells.Item(845, 13).Value = 0.0
$ws.Cells.Item(845, 14).Value = 0.0
$ws.Cells.Item(845, 15).Value = 0.0
$ws.Cells.Item(845, 16).Value = 4.86
$ws.Cells.Item(845, 17).Value = 24.36
$ws.Cells.Item(845, 18).Value = 4.3
$ws.Cells.Item(845, 19).Value = 31.0
$ws.Cells.Item(845, 20).Value = 5.0
$ws.Cells.Item(845, 21).Value = 18.0
$ws.Cells.Item(845, 22).Value = 5.0

# Row 846
$ws.Cells.Item(846, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(846, 6).Value = "left forward"
$ws.Cells.Item(846, 7).Value = "01:29:10"
$ws.Cells.Item(846, 8).Value = 8.15
$ws.Cells.Item(846, 9).Value = 0.27
$ws.Cells.Item(846, 10).Value = 7.87
$ws.Cells.Item(846, 11).Value = 0.22
$ws.Cells.Item(846, 12).Value = 0.06
$ws.Cells.Item(846, 13).Value = 0.0
$ws.Cells.Item(846, 14).Value = 0.0
$ws.Cells.Item(846, 15).Value = 2.0
$ws.Cells.Item(846, 16).Value = 4.24
$ws.Cells.Item(846, 17).Value = 25.28
$ws.Cells.Item(846, 18).Value = 4.75
$ws.Cells.Item(846, 19).Value = 51.0
$ws.Cells.Item(846, 20).Value = 14.0
$ws.Cells.Item(846, 21).Value = 39.0
$ws.Cells.Item(846, 22).Value = 8.0

# Row 847
$ws.Cells.Item(847, 5).Value = "Yoan Zouma"
$ws.Cells.Item(847, 6).Value = "center back"
$ws.Cells.Item(847, 7).Value = "01:25:58"
$ws.Cells.Item(847, 8).Value = 4.55
$ws.Cells.Item(847, 9).Value = 0.13
$ws.Cells.Item(847, 10).Value = 4.42
$ws.Cells.Item(847, 11).Value = 0.12
$ws.Cells.Item(847, 12).Value = 0.02
$ws.Cells.Item(847, 13).Value = 0.0
$ws.Cells.Item(847, 14).Value = 0.0
$ws.Cells.Item(847, 15).Value = 0.0
$ws.Cells.Item(847, 16).Value = 3.02
$ws.Cells.Item(847, 17).Value = 22.52
$ws.Cells.Item(847, 18).Value = 4.18
$ws.Cells.Item(847, 19).Value = 17.0
$ws.Cells.Item(847, 20).Value = 2.0
$ws.Cells.Item(847, 21).Value = 4.0
$ws.Cells.Item(847, 22).Value = 1.0

# Row 848
$ws.Cells.Item(848, 5).Value = "Omar Benyounes"
$ws.Cells.Item(848, 6).Value = "center midfield"
$ws.Cells.Item(848, 7).Value = "01:27:08"
$ws.Cells.Item(848, 8).Value = 7.05
$ws.Cells.Item(848, 9).Value = 0.43
$ws.Cells.Item(848, 10).Value = 6.61
$ws.Cells.Item(848, 11).Value = 0.36
$ws.Cells.Item(848, 12).Value = 0.08
$ws.Cells.Item(848, 13).Value = 0.0
$ws.Cells.Item(848, 14).Value = 0.0
$ws.Cells.Item(848, 15).Value = 0.0
$ws.Cells.Item(848, 16).Value = 4.8
$ws.Cells.Item(848, 17).Value = 24.68
$ws.Cells.Item(848, 18).Value = 4.32
$ws.Cells.Item(848, 19).Value = 32.0
$ws.Cells.Item(848, 20).Value = 4.0
$ws.Cells.Item(848, 21).Value = 27.0
$ws.Cells.Item(848, 22).Value = 3.0

# Row 849
$ws.Cells.Item(849, 5).Value = "Levy Ndoutoume"
$ws.Cells.Item(849, 6).Value = "left back"
$ws.Cells.Item(849, 7).Value = "01:29:28"
$ws.Cells.Item(849, 8).Value = 5.38
$ws.Cells.Item(849, 9).Value = 0.24
$ws.Cells.Item(849, 10).Value = 5.14
$ws.Cells.Item(849, 11).Value = 0.17
$ws.Cells.Item(849, 12).Value = 0.06
$ws.Cells.Item(849, 13).Value = 0.01
$ws.Cells.Item(849, 14).Value = 0.0
$ws.Cells.Item(849, 15).Value = 1.0
$ws.Cells.Item(849, 16).Value = 3.18
$ws.Cells.Item(849, 17).Value = 26.61
$ws.Cells.Item(849, 18).Value = 6.32
$ws.Cells.Item(849, 19).Value = 37.0
$ws.Cells.Item(849, 20).Value = 10.0
$ws.Cells.Item(849, 21).Value = 25.0
$ws.Cells.Item(849, 22).Value = 13.0

# Row 850
$ws.Cells.Item(850, 5).Value = "Jeremie Laurent"
$ws.Cells.Item(850, 6).Value = "left forward"
$ws.Cells.Item(850, 7).Value = "01:25:12"
$ws.Cells.Item(850, 8).Value = 6.05
$ws.Cells.Item(850, 9).Value = 0.38
$ws.Cells.Item(850, 10).Value = 5.66
$ws.Cells.Item(850, 11).Value = 0.33
$ws.Cells.Item(850, 12).Value = 0.06
$ws.Cells.Item(850, 13).Value = 0.0
$ws.Cells.Item(850, 14).Value = 0.0
$ws.Cells.Item(850, 15).Value = 0.0
$ws.Cells.Item(850, 16).Value = 4.17
$ws.Cells.Item(850, 17).Value = 24.75
$ws.Cells.Item(850, 18).Value = 5.09
$ws.Cells.Item(850, 19).Value = 39.0
$ws.Cells.Item(850, 20).Value = 13.0
$ws.Cells.Item(850, 21).Value = 39.0
$ws.Cells.Item(850, 22).Value = 10.0

# Row 851
$ws.Cells.Item(851, 5).Value = "Malik Boussaid"
$ws.Cells.Item(851, 6).Value = "right back"
$ws.Cells.Item(851, 7).Value = "01:31:53"
$ws.Cells.Item(851, 8).Value = 5.98
$ws.Cells.Item(851, 9).Value = 0.26
$ws.Cells.Item(851, 10).Value = 5.71
$ws.Cells.Item(851, 11).Value = 0.23
$ws.Cells.Item(851, 12).Value = 0.04
$ws.Cells.Item(851, 13).Value = 0.0
$ws.Cells.Item(851, 14).Value = 0.0
$ws.Cells.Item(851, 15).Value = 0.0
$ws.Cells.Item(851, 16).Value = 3.48
$ws.Cells.Item(851, 17).Value = 23.56
$ws.Cells.Item(851, 18).Value = 4.88
$ws.Cells.Item(851, 19).Value = 56.0
$ws.Cells.Item(851, 20).Value = 10.0
$ws.Cells.Item(851, 21).Value = 40.0
$ws.Cells.Item(851, 22).Value = 18.0

# Row 852
$ws.Cells.Item(852, 5).Value = "Ilyes Boughanmi"
$ws.Cells.Item(852, 6).Value = "center forward"
$ws.Cells.Item(852, 7).Value = "01:28:32"
$ws.Cells.Item(852, 8).Value = 5.78
$ws.Cells.Item(852, 9).Value = 0.14
$ws.Cells.Item(852, 10).Value = 5.63
$ws.Cells.Item(852, 11).Value = 0.14
$ws.Cells.Item(852, 12).Value = 0.01
$ws.Cells.Item(852, 13).Value = 0.0
$ws.Cells.Item(852, 14).Value = 0.0
$ws.Cells.Item(852, 15).Value = 1.0
$ws.Cells.Item(852, 16).Value = 3.29
$ws.Cells.Item(852, 17).Value = 25.45
$ws.Cells.Item(852, 18).Value = 5.02
$ws.Cells.Item(852, 19).Value = 34.0
$ws.Cells.Item(852, 20).Value = 8.0
$ws.Cells.Item(852, 21).Value = 40.0
$ws.Cells.Item(852, 22).Value = 10.0

# Row 853
$ws.Cells.Item(853, 5).Value = "Amir Etien"
$ws.Cells.Item(853, 6).Value = "right forward"
$ws.Cells.Item(853, 7).Value = "01:29:42"
$ws.Cells.Item(853, 8).Value = 5.22
$ws.Cells.Item(853, 9).Value = 0.32
$ws.Cells.Item(853, 10).Value = 4.89
$ws.Cells.Item(853, 11).Value = 0.22
$ws.Cells.Item(853, 12).Value = 0.08
$ws.Cells.Item(853, 13).Value = 0.02
$ws.Cells.Item(853, 14).Value = 0.0
$ws.Cells.Item(853, 15).Value = 3.0
$ws.Cells.Item(853, 16).Value = 2.82
$ws.Cells.Item(853, 17).Value = 30.36
$ws.Cells.Item(853, 18).Value = 5.85
$ws.Cells.Item(853, 19).Value = 43.0
$ws.Cells.Item(853, 20).Value = 13.0
$ws.Cells.Item(853, 21).Value = 23.0
$ws.Cells.Item(853, 22).Value = 12.0

# Move the view to mirror the author's final scroll/selection position.
$ws.Range("D860").Select()
